# Block_WBC.xlsx edit: "Try removing series with more than 20 percent missing observations"
#
# B1 = sheet1 (grid of rank positions per series + a 1..60 lookup table in K:M)
# B2 = sheet2 (small 2x2 summary table)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Remove series that have more than 20% missing observations from the
#    A2:I18 rank grid on sheet B1 (this also drops their rank numbers from
#    the count-lookup so the shared strings identifying dropped series fall
#    out of use and the remaining unique names stay in their original order).
# ---------------------------------------------------------------------------
$ws1.Range("G2").Value = 43
$ws1.Range("G3").Value = 45
$ws1.Range("H3").Value = 49
$ws1.Range("G4").ClearContents()
$ws1.Range("H4").Value = 50
$ws1.Range("A5").Value = 5
$ws1.Range("G5").ClearContents()
$ws1.Range("H5").Value = 51
$ws1.Range("A6").Value = 6
$ws1.Range("D6").Value = 26
$ws1.Range("G6").ClearContents()
$ws1.Range("H6").Value = 52
$ws1.Range("A7").ClearContents()
$ws1.Range("D7").ClearContents()
$ws1.Range("H7").Value = 53
$ws1.Range("A8").ClearContents()
$ws1.Range("F8").ClearContents()
$ws1.Range("H8").ClearContents()
$ws1.Range("C9").ClearContents()
$ws1.Range("C10").ClearContents()

# ---------------------------------------------------------------------------
# 2. Paste in the refreshed commodity-name list (column M, rows 1-56) -
#    these are the updated/renamed/reordered series labels used by the
#    COUNTIF lookup in column K. Writing top-to-bottom keeps the brand new
#    names appended to the shared-string table in this exact order.
# ---------------------------------------------------------------------------
$names = @(
    "Crude oil, Brent",
    "Crude oil, Dubai",
    "Crude oil, WTI",
    "Coal, Australian",
    "Natural gas, US",
    "Natural gas, Europe",
    "Liquefied natural gas, Japan",
    "Cocoa",
    "Coffee, Arabica",
    "Coffee, Robusta",
    "Tea, avg 3 auctions",
    "Coconut oil",
    "Copra",
    "Groundnuts",
    "Groundnut oil",
    "Palm oil",
    "Soybeans",
    "Soybean meal",
    "Rapeseed oil",
    "Sunflower oil",
    "Barley",
    "Maize",
    "Sorghum",
    "Rice, Thai",
    "Rice, Viet",
    "Wheat, US SRW",
    "Banana, US",
    "PSALM",
    "Fish meal",
    "Beef",
    "Meat, chicken",
    "Meat, sheep",
    "Shrimps, Mexican",
    "Sugar, world",
    "Logs, Cam",
    "Logs, Mal",
    "Sawnwood, Mal",
    "Plywood",
    "Woodpulp",
    "Tobacco, US",
    "Rubber, TSR20",
    "Phosphate rock",
    "DAP",
    "TSP",
    "Urea ",
    "Potassium chloride",
    "Aluminum",
    "Iron ore",
    "Copper",
    "Lead",
    "Tin",
    "Nickel",
    "Zinc",
    "Gold",
    "Platinum",
    "Silver"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 1
    $ws1.Range("M$row").Value = $names[$i]
}

# ---------------------------------------------------------------------------
# 3. Update the category header row (A1:I1) on B1 and the matching pair on
#    B2 to reference the same (re-pointed) shared strings.
# ---------------------------------------------------------------------------
$ws1.Range("A1").Value = "Energy"
$ws1.Range("B1").Value = "Beverages"
$ws1.Range("C1").Value = "Fats and Oils"
$ws1.Range("D1").Value = "Grains"
$ws1.Range("E1").Value = "Other food"
$ws1.Range("F1").Value = "Raw materials"
$ws1.Range("G1").Value = "Fertilizers"
$ws1.Range("H1").Value = "Metals"
$ws1.Range("I1").Value = "Precious metals"

$ws2.Range("A1").Value = "Beverages"
$ws2.Range("B1").Value = "Fats and Oils"

# ---------------------------------------------------------------------------
# 4. Restore the view: B1 becomes the active/selected tab with G11 selected,
#    B2 loses tab-selection and keeps its own C3 selection.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("G11").Select()
$ws2.Range("C3").Select()
